# Remove pushofbiz / cabitest5 credentials from the workbook and replace
# with the sanitized cabitest3 data, per commit:
# "removed pushofbiz credentials from excel data file"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G2: backoffice control URL - test19 -> test3
$ws.Range("G2").Value2 = "https://test3.cliotest.com/backoffice/control/main"

# H2: vhost id - cabitest5 -> cabitest3
$ws.Range("H2").Value2 = "cabitest3"

# A3: show microsite link - cabitest5/104516894 -> cabitest3/104526592
$ws.Range("A3").Value2 = "https://mirandakate.cabitest3.com/show-microsite/104526592/"

# B6: ofbiz user - sshinde -> abcd
$ws.Range("B6").Value2 = "abcd"

# C6: ofbiz password - mask the visible value, but keep the original
# mailto hyperlink friendly text recorded as the "display" attribute,
# same as Excel does when a hyperlink's display text no longer matches
# its address text.
foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$C$6') {
        $h.TextToDisplay = "C@bi`$ush5"
    }
}
$ws.Range("C6").Value2 = "******"

# Update the active selection to B6 (matches the saved sheetView state)
$ws.Range("B6").Select()
